$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.905.63'
$ws.Range("E2").Value = '  -0.93%  '
$ws.Range("D3").Value = '2.901.66'
$ws.Range("E3").Value = '  -1.15%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '569.32'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -3.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.70'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.32%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  -1.02%  '
$ws.Range("D9").Value = '2.900.96'
$ws.Range("E9").Value = '  -1.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.95'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.75%  '
$ws.Range("E11").Value = '  -2.78%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.429'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.57%  '
$ws.Range("E13").Value = '  -1.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.15'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.36%  '
$ws.Range("E15").Value = '  -0.20%  '
$ws.Range("D16").Value = '3.382.41'
$ws.Range("E16").Value = '  -1.16%  '
$ws.Range("D17").Value = '61.846.39'
$ws.Range("E17").Value = '  -0.99%  '
$ws.Range("D18").Value = '2.898.11'
$ws.Range("E18").Value = '  -1.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.51'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.79%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '429.83'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.84%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.91'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -4.01%  '
$ws.Range("E22").Value = '  -1.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.86'
$ws.Range("D23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.01'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.47%  '
$ws.Range("E25").Value = '  +0.67%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.11'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -8.07%  '
$ws.Range("E27").Value = '  +0.15%  '
$ws.Range("E28").Value = '  -3.02%  '
$ws.Range("E29").Value = '  +7.46%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.01'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.80%  '
$ws.Range("E31").Value = '  -2.30%  '
$ws.Range("E32").Value = '  -6.18%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("E34").Value = '  -3.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '25.62'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.64%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.955'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -3.19%  '
$ws.Range("E37").Value = '  -3.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '48.84'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.50%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.84'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -5.61%  '
$ws.Range("E40").Value = '  -4.51%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.114'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.28%  '
$ws.Range("E42").Value = '  -2.44%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.74'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +4.17%  '
$ws.Range("E44").Value = '  -2.45%  '
$ws.Range("D45").Value = '2.704.75'
$ws.Range("E45").Value = '  +0.40%  '
$ws.Range("E46").Value = '  -0.85%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '131.76'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.60%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '346.17'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.59%  '
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("E50").Value = '  -0.91%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.56'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -4.07%  '
